# Add data for 2022-07-14
# This updates the "through" date in the sheet title / header (July 05 -> July 06)
# and adds newly-recorded carjacking counts (for 2022-07-14 data refresh) to the
# "July" columns of several prior years across a handful of neighborhoods.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-07-06"

# Update the column header text that mirrors the "through" date
$ws.Range("B1").Value2 = "July 2022 (through July 06)"

# --- Data updates -----------------------------------------------------

# Row 2: Austin
$ws.Range("I2").Value2 = 3
$ws.Range("AD2").Value2 = 2

# Row 3: Englewood
$ws.Range("P3").Value2 = 1

# Row 5: Garfield Park
$ws.Range("AR5").Value2 = 1

# Row 7: Roseland
$ws.Range("AD7").Value2 = 2

# Row 8: North Lawndale
$ws.Range("P8").Value2 = 3

# Row 15: Douglas
$ws.Range("P15").Value2 = 1

# Row 20: Kenwood
$ws.Range("I20").Value2 = 1

# Row 26: Little Village
$ws.Range("B26").Value2 = 1

# Row 29: Humboldt Park
$ws.Range("AD29").Value2 = 1

# Row 38: West Town
$ws.Range("AR38").Value2 = 1

# Row 49: Grand Boulevard
$ws.Range("I49").Value2 = 2

# Row 52: Chatham
$ws.Range("I52").Value2 = 3
$ws.Range("P52").Value2 = 3

# Row 53: Calumet Heights
$ws.Range("AR53").Value2 = 1

# Row 70: Edgewater
$ws.Range("AD70").Value2 = 1

# Row 96: Wrigleyville
$ws.Range("I96").Value2 = 3
